#
# click-widget-model.pptx edit:
#  1. Refresh the cached "datetimeFigureOut" date field (master + every
#     slide layout) from 18/01/2021 to 20/01/2021.
#  2. Swap the X positions of the two vertical connector arrows
#     ("Connettore 2 28" / "Connettore 2 30") on slide 1.
#

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh on the slide master and every layout.
# ---------------------------------------------------------------------

$newDate = "20/01/2021"

function Update-DateShapes($shapeColl, $newText) {
    for ($idxShp = 1; $idxShp -le $shapeColl.Count; $idxShp++) {
        $shp = $shapeColl.Item($idxShp)
        if ($shp.HasTextFrame) {
            $isDatePh = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
            if ($isDatePh) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($idxLyt = 1; $idxLyt -le $layouts.Count; $idxLyt++) {
    $lyt = $layouts.Item($idxLyt)
    Update-DateShapes $lyt.Shapes $newDate
}

# ---------------------------------------------------------------------
# 2) Re-position the two connectors on slide 1. Shape.Left is a Single
#    (float32) under the hood, so the literal point values below are
#    chosen so they round-trip to the exact target EMU offsets
#    (1619672 and 899592 respectively).
# ---------------------------------------------------------------------

$slide = $p.Slides.Item(1)

$connA = $slide.Shapes.Item("Connettore 2 28")
$connA.Left = 127.5333

$connB = $slide.Shapes.Item("Connettore 2 30")
$connB.Left = 70.83406
